$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = "705"
$ws.Range("P3").Value = "954155"

$ws.Range("O4").Value = "554"
$ws.Range("P4").Value = "908831"

$ws.Range("O5").Value = "538"
$ws.Range("P5").Value = "722227"

$ws.Range("O6").Value = "529"
$ws.Range("P6").Value = "530561"

$ws.Range("O7").Value = "557"
$ws.Range("P7").Value = "789674"

$ws.Range("O8").Value = "853"
$ws.Range("P8").Value = "885596"

$ws.Range("O9").Value = "943"
$ws.Range("P9").Value = "750585"

$ws.Range("O10").Value = "387"
$ws.Range("P10").Value = "814129"

$ws.Range("O11").Value = "513"
$ws.Range("P11").Value = "822893"

$ws.Range("O12").Value = "536"
$ws.Range("P12").Value = "1024.9k"

$ws.Range("O13").Value = "581"
$ws.Range("P13").Value = "1371.2k"

$ws.Range("O14").Value = "819"
$ws.Range("P14").Value = "658730"
